# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for "Membrillo" (Vega Modelo de Temuco)
# above the current row 178, shifting the existing rows 178-212 down to
# 179-213 (dimension grows from A1:T212 to A1:T213).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 178:212 down one row, creating a blank row 178.
$ws.Rows(178).Insert()

# Populate the newly inserted row 178 with the new weekly record.
$ws.Cells.Item(178, 1).Value  = 10
$ws.Cells.Item(178, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(178, 3).Value  = "La Araucanía"
$ws.Cells.Item(178, 4).Value  = 44785
$ws.Cells.Item(178, 5).Value  = 9
$ws.Cells.Item(178, 6).Value  = "Fruta"
$ws.Cells.Item(178, 7).Value  = 100104
$ws.Cells.Item(178, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(178, 9).Value  = 100104003
$ws.Cells.Item(178, 10).Value = "Membrillo"
$ws.Cells.Item(178, 11).Value = "Champion"
$ws.Cells.Item(178, 12).Value = "Primera"
$ws.Cells.Item(178, 13).Value = 155
$ws.Cells.Item(178, 14).Value = 10000
$ws.Cells.Item(178, 15).Value = 10000
$ws.Cells.Item(178, 16).Value = 10000
$ws.Cells.Item(178, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(178, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(178, 19).Value = 556
$ws.Cells.Item(178, 20).Value = 18
